$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: find the index (1-based) of the first paragraph whose text contains
# a given substring.
# ---------------------------------------------------------------------------
function Find-ParagraphIndex($needle) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text -like "*$needle*") {
            return $i
        }
    }
    return -1
}

function New-OpenXmlPackage($bodyXml) {
    return "<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'><pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'><pkg:xmlData><w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:body>$bodyXml</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"
}

# ---------------------------------------------------------------------------
# 1) Insert a new paragraph (plus its own trailing blank separator paragraph)
#    right before the "Very rarely will you be creating..." paragraph,
#    mentioning null-value imputation / interpolation techniques.
# ---------------------------------------------------------------------------
$veryRarelyIdx = Find-ParagraphIndex("Very rarely will you be creating")
if ($veryRarelyIdx -ge 1) {
    $target = $d.Paragraphs.Item($veryRarelyIdx)
    $insPoint = $d.Range($target.Range.Start, $target.Range.Start)

    $newParaText = "Another important step is filling in null values. Ideally your client or vendor should give you a full dataset, but in the real world you will have null values and it’s just unavoidable. Easy ways to resolve null values is to replace them with a mean or median value for that feature, but you can also try out more advanced techniques such as Vandermonde interpolation, cubic spline, Newtonian interpolation, Lagrange interpolation, etc."

    $bodyXml = "<w:p><w:r><w:t>PLACEHOLDER_TEXT</w:t></w:r></w:p><w:p/>"
    $bodyXml = $bodyXml.Replace("PLACEHOLDER_TEXT", $newParaText)
    $insPoint.InsertXML((New-OpenXmlPackage $bodyXml))
}

# ---------------------------------------------------------------------------
# 2) Move the <w:lastRenderedPageBreak/> marker: it currently sits on the
#    "Reflect on the interrelationship..." paragraph, but it should instead
#    be on the "I think in general it's useful..." paragraph.
# ---------------------------------------------------------------------------

# -- 2a) Add lastRenderedPageBreak to the "I think in general" run.
$iThinkIdx = Find-ParagraphIndex("I think in general it")
if ($iThinkIdx -ge 1) {
    $p = $d.Paragraphs.Item($iThinkIdx)
    $r = $p.Range
    $r.MoveEnd(1, -1) | Out-Null   # exclude the paragraph mark (keeps pPr)
    $text = $r.Text
    $start = $r.Start
    $r.Delete()

    $insertionPoint = $d.Range($start, $start)
    $bodyXml = "<w:p><w:r><w:lastRenderedPageBreak/><w:t>PLACEHOLDER_TEXT</w:t></w:r></w:p>"
    $bodyXml = $bodyXml.Replace("PLACEHOLDER_TEXT", $text)
    $insertionPoint.InsertXML((New-OpenXmlPackage $bodyXml))
}

# -- 2b) Remove lastRenderedPageBreak from the "Reflect on the
#        interrelationship..." run.
$reflectIdx = Find-ParagraphIndex("Reflect on the interrelationship")
if ($reflectIdx -ge 1) {
    $p = $d.Paragraphs.Item($reflectIdx)
    $r = $p.Range
    $r.MoveEnd(1, -1) | Out-Null   # exclude the paragraph mark (keeps pPr/numPr)
    $text = $r.Text
    $r.Delete()
    $p.Range.InsertBefore($text)
}
